$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 76
$ws.Cells.Item(76, 8).Value = 3607.8333
$ws.Cells.Item(76, 10).Value = 3327
$ws.Cells.Item(76, 12).Value = 3327
$ws.Cells.Item(76, 14).Value = -3957
# row 79
$ws.Cells.Item(79, 8).Value = 3607.8333
$ws.Cells.Item(79, 10).Value = 3327
$ws.Cells.Item(79, 12).Value = 3327
$ws.Cells.Item(79, 14).Value = -5511
# row 86
$ws.Cells.Item(86, 8).Value = 8251.846
$ws.Cells.Item(86, 9).Value = 3765.5
$ws.Cells.Item(86, 11).Value = 3765.5
$ws.Cells.Item(86, 13).Value = -2642.5
# row 89
$ws.Cells.Item(89, 8).Value = 8251.846
$ws.Cells.Item(89, 9).Value = 3765.5
$ws.Cells.Item(89, 11).Value = 18827.5
$ws.Cells.Item(89, 13).Value = -13211.5
# row 106
$ws.Cells.Item(106, 8).Value = 6201.3335
$ws.Cells.Item(106, 9).Value = 4302
$ws.Cells.Item(106, 11).Value = 4302
$ws.Cells.Item(106, 13).Value = -3671
# row 132
$ws.Cells.Item(132, 8).Value = 1711
$ws.Cells.Item(132, 9).Value = 1712.375
$ws.Cells.Item(132, 11).Value = 5137.125
$ws.Cells.Item(132, 13).Value = -2607.125
# row 137
$ws.Cells.Item(137, 8).Value = 1693.25
$ws.Cells.Item(137, 9).Value = 1386.5
$ws.Cells.Item(137, 11).Value = 4159.5
$ws.Cells.Item(137, 13).Value = -1609.5
# row 138
$ws.Cells.Item(138, 8).Value = 1976.3636
$ws.Cells.Item(138, 9).Value = 540
$ws.Cells.Item(138, 10).Value = 3700
$ws.Cells.Item(138, 11).Value = 1620
$ws.Cells.Item(138, 12).Value = 11100
$ws.Cells.Item(138, 13).Value = 3520
$ws.Cells.Item(138, 14).Value = -21380
# row 141
$ws.Cells.Item(141, 8).Value = 2313.7307
$ws.Cells.Item(141, 9).Value = 2154.652
$ws.Cells.Item(141, 11).Value = 6463.956
$ws.Cells.Item(141, 13).Value = -1283.956

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 38
$ws.Cells.Item(38, 8).Value = 3999
$ws.Cells.Item(38, 9).Value = 3999
$ws.Cells.Item(38, 11).Value = 3999
$ws.Cells.Item(38, 13).Value = -3532
# row 61
$ws.Cells.Item(61, 8).Value = 4530
$ws.Cells.Item(61, 9).Value = 4295
$ws.Cells.Item(61, 11).Value = 4295
$ws.Cells.Item(61, 13).Value = -4083
# row 74
$ws.Cells.Item(74, 8).Value = 6700
# row 76
$ws.Cells.Item(76, 8).Value = 45000
$ws.Cells.Item(76, 10).Value = 45000
$ws.Cells.Item(76, 12).Value = 45000
$ws.Cells.Item(76, 14).Value = -45676
# row 77
$ws.Cells.Item(77, 8).Value = 6700
# row 79
$ws.Cells.Item(79, 8).Value = 45000
$ws.Cells.Item(79, 10).Value = 45000
$ws.Cells.Item(79, 12).Value = 45000
$ws.Cells.Item(79, 14).Value = -47340
# row 97
$ws.Cells.Item(97, 8).Value = 930.44446
$ws.Cells.Item(97, 9).Value = 930.44446
$ws.Cells.Item(97, 11).Value = 930.44446
$ws.Cells.Item(97, 13).Value = -434.44446
# row 122
$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 13).ClearContents()
# row 136
$ws.Cells.Item(136, 8).Value = 4530
$ws.Cells.Item(136, 9).Value = 4295
$ws.Cells.Item(136, 11).Value = 12885
$ws.Cells.Item(136, 13).Value = -10335

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 5
$ws.Cells.Item(5, 8).Value = 7170.6
$ws.Cells.Item(5, 9).Value = 849
$ws.Cells.Item(5, 10).Value = 8751
$ws.Cells.Item(5, 11).Value = 849
$ws.Cells.Item(5, 12).Value = 8751
$ws.Cells.Item(5, 13).Value = -736
$ws.Cells.Item(5, 14).Value = -8977
# row 22
$ws.Cells.Item(22, 8).Value = 399.15384
$ws.Cells.Item(22, 9).Value = 399.15384
$ws.Cells.Item(22, 11).Value = 399.15384
$ws.Cells.Item(22, 13).Value = -226.15384
# row 112
$ws.Cells.Item(112, 8).Value = 24489.334
$ws.Cells.Item(112, 10).Value = 24489.334
$ws.Cells.Item(112, 12).Value = 24489.334
$ws.Cells.Item(112, 14).Value = -27443.334
# row 134
$ws.Cells.Item(134, 8).Value = 7278
$ws.Cells.Item(134, 9).Value = 8205.833000000001
$ws.Cells.Item(134, 11).Value = 24617.499
$ws.Cells.Item(134, 13).Value = -22082.499

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 132
$ws.Cells.Item(132, 8).Value = 1493
$ws.Cells.Item(132, 9).Value = 1493
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 4479
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -1949
$ws.Cells.Item(132, 14).ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 26
$ws.Cells.Item(26, 8).Value = 153642.86
$ws.Cells.Item(26, 9).Value = 335166.34
$ws.Cells.Item(26, 11).Value = 1005499.02
$ws.Cells.Item(26, 13).Value = -1005211.02
# row 80
$ws.Cells.Item(80, 8).Value = 7955.6665
$ws.Cells.Item(80, 9).Value = 1798
$ws.Cells.Item(80, 10).Value = 9187.200000000001
$ws.Cells.Item(80, 11).Value = 5394
$ws.Cells.Item(80, 12).Value = 27561.6
$ws.Cells.Item(80, 13).Value = -4458
$ws.Cells.Item(80, 14).Value = -29433.6
# row 83
$ws.Cells.Item(83, 8).Value = 7955.6665
$ws.Cells.Item(83, 9).Value = 1798
$ws.Cells.Item(83, 10).Value = 9187.200000000001
$ws.Cells.Item(83, 11).Value = 16182
$ws.Cells.Item(83, 12).Value = 82684.8
$ws.Cells.Item(83, 13).Value = -11502
$ws.Cells.Item(83, 14).Value = -92044.8
# row 138
$ws.Cells.Item(138, 8).Value = 5140.091
$ws.Cells.Item(138, 10).Value = 4706.375
$ws.Cells.Item(138, 12).Value = 14119.125
$ws.Cells.Item(138, 14).Value = -24399.125

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 109
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 14).ClearContents()
# row 122
$ws.Cells.Item(122, 8).Value = 10419731
$ws.Cells.Item(122, 9).Value = 13890529
$ws.Cells.Item(122, 11).Value = 41671587
$ws.Cells.Item(122, 13).Value = -41669137
# row 132
$ws.Cells.Item(132, 8).Value = 3704.5
$ws.Cells.Item(132, 9).Value = 3704.5
$ws.Cells.Item(132, 11).Value = 11113.5
$ws.Cells.Item(132, 13).Value = -8583.5
# row 135
$ws.Cells.Item(135, 10).Value = 48779.2
$ws.Cells.Item(135, 12).Value = 48779.2
$ws.Cells.Item(135, 14).Value = -58919.2

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Cells.Item(22, 8).Value = 2806
$ws.Cells.Item(22, 9).Value = 2806
$ws.Cells.Item(22, 11).Value = 2806
$ws.Cells.Item(22, 13).Value = -2511
# row 23
$ws.Cells.Item(23, 8).Value = 25000
$ws.Cells.Item(23, 10).Value = 25000
$ws.Cells.Item(23, 12).Value = 25000
$ws.Cells.Item(23, 14).Value = -25460
# row 27
$ws.Cells.Item(27, 8).Value = 2806
$ws.Cells.Item(27, 9).Value = 2806
$ws.Cells.Item(27, 11).Value = 2806
$ws.Cells.Item(27, 13).Value = -2699
# row 30
$ws.Cells.Item(30, 8).Value = 5488.25
$ws.Cells.Item(30, 9).Value = 645
$ws.Cells.Item(30, 10).Value = 20018
$ws.Cells.Item(30, 11).Value = 645
$ws.Cells.Item(30, 12).Value = 20018
$ws.Cells.Item(30, 13).Value = -537
$ws.Cells.Item(30, 14).Value = -20234
# row 40
$ws.Cells.Item(40, 8).Value = 2999
$ws.Cells.Item(40, 9).Value = 2999
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 2999
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = -2863
$ws.Cells.Item(40, 14).ClearContents()
# row 68
$ws.Cells.Item(68, 8).Value = 1723
$ws.Cells.Item(68, 9).Value = 1723
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 1723
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).Value = -974
$ws.Cells.Item(68, 14).ClearContents()
# row 71
$ws.Cells.Item(71, 8).Value = 1723
$ws.Cells.Item(71, 9).Value = 1723
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 8615
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).Value = -4871
$ws.Cells.Item(71, 14).ClearContents()
# row 121
$ws.Cells.Item(121, 8).Value = 149995
$ws.Cells.Item(121, 10).Value = 149995
$ws.Cells.Item(121, 12).Value = 149995
$ws.Cells.Item(121, 14).Value = -153489
# row 132
$ws.Cells.Item(132, 8).Value = 4280.4736
$ws.Cells.Item(132, 9).Value = 3171.4614
$ws.Cells.Item(132, 10).Value = 6683.3335
$ws.Cells.Item(132, 11).Value = 9514.3842
$ws.Cells.Item(132, 12).Value = 20050.0005
$ws.Cells.Item(132, 13).Value = -6984.3842
$ws.Cells.Item(132, 14).Value = -25110.0005

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 69
$ws.Cells.Item(69, 8).Value = 23114.2
$ws.Cells.Item(69, 10).Value = 23114.2
$ws.Cells.Item(69, 12).Value = 23114.2
$ws.Cells.Item(69, 14).Value = -24612.2
# row 72
$ws.Cells.Item(72, 8).Value = 23114.2
$ws.Cells.Item(72, 10).Value = 23114.2
$ws.Cells.Item(72, 12).Value = 69342.60000000001
$ws.Cells.Item(72, 14).Value = -76830.60000000001
# row 112
$ws.Cells.Item(112, 8).Value = 50000
$ws.Cells.Item(112, 10).Value = 50000
$ws.Cells.Item(112, 12).Value = 50000
$ws.Cells.Item(112, 14).Value = -52954
# row 121
$ws.Cells.Item(121, 8).Value = 124997.5
$ws.Cells.Item(121, 10).Value = 124997.5
$ws.Cells.Item(121, 12).Value = 124997.5
$ws.Cells.Item(121, 14).Value = -128491.5
# row 122
$ws.Cells.Item(122, 8).Value = 500
$ws.Cells.Item(122, 9).Value = 500
$ws.Cells.Item(122, 11).Value = 1500
$ws.Cells.Item(122, 13).Value = 950
# row 132
$ws.Cells.Item(132, 8).Value = 1743
$ws.Cells.Item(132, 9).Value = 1743
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 5229
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -2699
$ws.Cells.Item(132, 14).ClearContents()
